$d = $word.ActiveDocument

function Get-LastMatchRange($doc, $searchText) {
    $lastMatch = $null
    $rng = $doc.Content
    $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    while ($rng.Find.Found) {
        $lastMatch = $doc.Range($rng.Start, $rng.End)
        $rng.Collapse(0)
        $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    }
    return $lastMatch
}

# --- "Eigene Leistung:" (last occurrence, in the "Abschlusswoche" section) ---
$eigeneLeistung = Get-LastMatchRange $d "Eigene Leistung:"
$insertPos = $eigeneLeistung.End
$r1 = $d.Range($insertPos, $insertPos)
$r1.InsertAfter(" ")
$insertPos2 = $insertPos + 1
$r2 = $d.Range($insertPos2, $insertPos2)
$r2.InsertAfter("Gravitation implementiert, Effizientere Berechnung (Quake III)")

# --- "Gruppenleistung:" (last occurrence, in the "Abschlusswoche" section) ---
$gruppenleistung = Get-LastMatchRange $d "Gruppenleistung:"
$insertPos = $gruppenleistung.End
$r3 = $d.Range($insertPos, $insertPos)
$r3.InsertAfter(" ")
$insertPos2 = $insertPos + 1
$r4 = $d.Range($insertPos2, $insertPos2)
$r4.InsertAfter("Gravitation implementiert, Effizientere Berechnung (Quake III), Präsentation")

# --- "Eigener Zeitaufwand:" (last occurrence, in the "Abschlusswoche" section) ---
$zeitaufwand = Get-LastMatchRange $d "Eigener Zeitaufwand:"
$insertPos = $zeitaufwand.End
$r5 = $d.Range($insertPos, $insertPos)
$r5.InsertAfter(" 2 Schulstunden zuhause, 4 Schulstunden in der Schule")
